# Adds "Query By" (C Madhu), "Query Added Date" (2020-06-26) and
# "Query Assigned To" (Ashok) data to the first 4 data rows (rows 2-5),
# columns H, I and K respectively. Also widens those columns to fit the
# new content and moves the active selection to K5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - 5: H = "C Madhu", I = date 2020-06-26, K = "Ashok"
foreach ($r in 2..5) {
    $ws.Cells.Item($r, 8).Value = "C Madhu"
    $ws.Cells.Item($r, 9).Value = 44008
    $ws.Cells.Item($r, 11).Value = "Ashok"
}

# Apply the date number format to I2 first, then re-apply the resulting
# (canonical) format string to the whole I2:I5 block in a single call so
# every cell converges on the same shared style (avoids minting a new,
# functionally-identical style per cell).
$ws.Range("I2").NumberFormat = "mm-dd-yy"
$dateFmt = $ws.Range("I2").NumberFormat
$ws.Range("I2:I5").NumberFormat = $dateFmt

# Resize columns to fit the newly entered content.
$ws.Columns.Item(8).ColumnWidth = 23.28515625
$ws.Columns.Item(10).ColumnWidth = 24.85546875
$ws.Columns.Item(11).ColumnWidth = 24.28515625

# Move the active selection to K5, matching the author's final cursor
# position after entering the data.
$ws.Range("K5").Select()
